# Auto-generated: update crypto price/volume columns (D, E) per upstream refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.122.96"
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").Value = "1.748.75"
$ws.Range("E3").Value = "  +0.08%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'236.93"
$ws.Range("E5").Value = "  +1.17%  "
$ws.Range("D6").Value = "'0.5609"
$ws.Range("E6").Value = "  +8.01%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'0.2834"
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("D9").Value = "'0.06186"
$ws.Range("E9").Value = "  +0.87%  "
$ws.Range("D10").Value = "1.755.09"
$ws.Range("E10").Value = "  +0.06%  "
$ws.Range("D11").Value = "'0.07192"
$ws.Range("E11").Value = "  +2.51%  "
$ws.Range("E12").Value = "  +0.52%  "
$ws.Range("D13").Value = "'0.6529"
$ws.Range("E13").Value = "  +1.57%  "
$ws.Range("D14").Value = "'4.640"
$ws.Range("E14").Value = "  +2.50%  "
$ws.Range("D15").Value = "'78.30"
$ws.Range("E15").Value = "  +0.97%  "
$ws.Range("E16").Value = "  -0.07%  "
$ws.Range("E17").Value = "  +0.08%  "
$ws.Range("D18").Value = "26.015.44"
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("D19").Value = "'11.80"
$ws.Range("E19").Value = "  +2.40%  "
$ws.Range("D20").Value = "'0.000006787"
$ws.Range("E20").Value = "  +2.52%  "
$ws.Range("D21").Value = "1.978.11"
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("D22").Value = "'4.366"
$ws.Range("E22").Value = "  +5.20%  "
$ws.Range("E23").Value = "  +1.11%  "
$ws.Range("D24").Value = "'5.302"
$ws.Range("E24").Value = "  +3.18%  "
$ws.Range("E25").Value = "  +0.37%  "
$ws.Range("D26").Value = "'1.526"
$ws.Range("E26").Value = "  +1.93%  "
$ws.Range("D27").Value = "'15.35"
$ws.Range("E27").Value = "  +1.46%  "
$ws.Range("E28").Value = "  -0.91%  "
$ws.Range("D29").Value = "'105.43"
$ws.Range("E29").Value = "  +2.60%  "
$ws.Range("D30").Value = "'0.08503"
$ws.Range("E30").Value = "  +2.65%  "
$ws.Range("D31").Value = "'3.794"
$ws.Range("E31").Value = "  +3.53%  "
$ws.Range("D32").Value = "'3.639"
$ws.Range("E32").Value = "  +5.65%  "
$ws.Range("D33").Value = "'0.04654"
$ws.Range("E33").Value = "  +3.92%  "
$ws.Range("D34").Value = "'2.650"
$ws.Range("E34").Value = "  +1.48%  "
$ws.Range("E35").Value = "  +1.72%  "
$ws.Range("D36").Value = "'0.6295"
$ws.Range("E36").Value = "  +1.95%  "
$ws.Range("D37").Value = "'2.706"
$ws.Range("E37").Value = "  +1.25%  "
$ws.Range("D38").Value = "'0.01610"
$ws.Range("E38").Value = "  +0.94%  "
$ws.Range("D39").Value = "'1.966"
$ws.Range("E39").Value = "  +1.38%  "
$ws.Range("D41").Value = "'100.92"
$ws.Range("E41").Value = "  +0.59%  "
$ws.Range("D42").Value = "'0.3937"
$ws.Range("E42").Value = "  +1.89%  "
$ws.Range("D43").Value = "'0.7479"
$ws.Range("E43").Value = "  -0.08%  "
$ws.Range("D44").Value = "'5.068"
$ws.Range("E44").Value = "  -0.12%  "
$ws.Range("D45").Value = "'0.1150"
$ws.Range("E45").Value = "  +1.96%  "
$ws.Range("D46").Value = "'6.353"
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("D47").Value = "'0.05335"
$ws.Range("E47").Value = "  -2.25%  "
$ws.Range("D48").Value = "'54.84"
$ws.Range("E48").Value = "  +3.34%  "
$ws.Range("D49").Value = "'30.84"
$ws.Range("E49").Value = "  +2.42%  "
$ws.Range("D50").Value = "'0.3494"
$ws.Range("E50").Value = "  +1.92%  "
$ws.Range("D51").Value = "'7.579"
$ws.Range("E51").Value = "  -0.61%  "
